# Correção nos dados e inicio da analise PNAD 2009
#
# The two "section header" rows (row 5 "situação do domicílio" and row 8
# "grandes regiões e unidades da federação") were empty placeholder rows
# with no B:H data. They are removed entirely and the data beneath them
# (urbana/rural, norte/rondônia/...) shifts up to take their place,
# acquiring the statistics that used to sit one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the row-8 reference is still
# valid when we delete it.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()

# Column F of the header row 2 used to read "total" while column B read
# "unnamed: 1_level_1" (an artifact label for the unlabeled merged header
# cell). Both now read "total".
$ws.Range("F2").ClearContents()
$ws.Range("B2").Value2 = "total"
$ws.Range("F2").Value2 = "total"
